$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.893344666666666
$ws.Range("H2").Value = 5.680033999999999
$ws.Range("I2").Value = 0.05525983881677096
$ws.Range("J2").Value = 0.05525983881677096
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3389413333333333
$ws.Range("N2").Value = 1.016824
$ws.Range("Q2").Value = 0.6417327657795554
$ws.Range("R2").Value = 5.775594892015999
$ws.Range("S2").Value = 0.05525983881677096
$ws.Range("T2").Value = 0.05525983881677096

# Row 3
$ws.Range("G3").Value = 4.159773333333334
$ws.Range("I3").Value = 0.1214086415227279
$ws.Range("J3").Value = 0.1214086415227279
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3389413333333333
$ws.Range("N3").Value = 1.016824
$ws.Range("Q3").Value = 1.409919119964445
$ws.Range("R3").Value = 12.68927207968
$ws.Range("S3").Value = 0.1214086415227279
$ws.Range("T3").Value = 0.1214086415227279

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.240212333333333
$ws.Range("H4").Value = 3.720637
$ws.Range("I4").Value = 0.03619728348733726
$ws.Range("J4").Value = 0.03619728348733727
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3389413333333333
$ws.Range("N4").Value = 1.016824
$ws.Range("Q4").Value = 0.4203592218764444
$ws.Range("R4").Value = 3.783232996888
$ws.Range("S4").Value = 0.03619728348733726
$ws.Range("T4").Value = 0.03619728348733727

# Row 5
$ws.Range("G5").Value = 26.96925
$ws.Range("H5").Value = 80.90774999999999
$ws.Range("I5").Value = 0.7871342361731639
$ws.Range("J5").Value = 0.7871342361731638
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3389413333333333
$ws.Range("N5").Value = 1.016824
$ws.Range("Q5").Value = 9.140993554
$ws.Range("R5").Value = 82.26894198599999
$ws.Range("S5").Value = 0.7871342361731639
$ws.Range("T5").Value = 0.7871342361731638
